$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first block of the diff
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1049
$ws1.Range("F3").Value = 410
$ws1.Range("F4").Value = 3103
$ws1.Range("F5").Value = 77
$ws1.Range("F6").Value = 636

# Sheet "全部类型" (All types) - second block of the diff
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1049
$ws4.Range("F4").Value = 410
$ws4.Range("F5").Value = 3103
$ws4.Range("F6").Value = 77
$ws4.Range("F7").Value = 636
